$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Five new accelerometer readings get inserted right after the header row,
# pushing the existing data rows down by five.
$newData = @(
    @(1.809551620483398, -2.747114396095276, -2.346350741386413),
    @(1.535980415344238, -2.904552030563355, -2.271616220474243),
    @(1.500288486480713, -2.749492883682251, -1.919559156894683),
    @(1.063152790069579, -2.663525581359863, -1.874630331993103),
    @(1.008758783340455, -2.917640089988708, -2.169865667819977)
)

$ws.Rows("2:6").Insert()

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newData[$i][0]
    $ws.Cells.Item($r, 2).Value = $newData[$i][1]
    $ws.Cells.Item($r, 3).Value = $newData[$i][2]
}

# Insert() copies the formatting of the row above (the bold header row) onto
# the freshly-inserted rows; strip that back off so the new data rows look
# like the rest of the plain data rows.
$ws.Range("A2:C6").ClearFormats()

# The old last six data rows (previously rows 17-22, now shifted to 22-27)
# are dropped so that the sheet ends up with 20 data rows (A1:C21 overall).
$ws.Rows("22:27").Delete()
